$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "White River Lower 08"
$ws.Range("B2").Value = "Wenatchee"
$ws.Range("C2").Value = "Lower White River"

$ws.Range("K2").Value = 3
$ws.Range("N2").Value = 1
$ws.Range("O2").Value = 1
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 4
$ws.Range("S2").Value = 5
$ws.Range("T2").Value = 34
$ws.Range("U2").Value = 0.7555555555555555
